$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("B248")
$fc = $rng.FormatConditions.Item(1)
$fc.Delete()

$new = $rng.FormatConditions.AddUniqueValues()
$new.DupeUnique = 1
$new.Borders.Item(7).LineStyle = 1
$new.Borders.Item(7).Color = 255
$new.Borders.Item(10).LineStyle = 1
$new.Borders.Item(10).Color = 255
$new.Borders.Item(8).LineStyle = 1
$new.Borders.Item(8).Color = 255
$new.Borders.Item(9).LineStyle = 1
$new.Borders.Item(9).Color = 255
$new.Priority = 1
